$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.760.67"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.448.81"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.60"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.75"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "2.444.49"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.11"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("E15").Value = "  -4.85%  "
$ws.Range("D16").Value = "2.890.99"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "62.384.05"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "2.439.71"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.29"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.60"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.12"
$ws.Range("E23").Value = "  +11.38%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.36"
$ws.Range("E25").Value = "  -3.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "627.94"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000100"
$ws.Range("E28").Value = "  -5.21%  "
$ws.Range("D29").Value = "2.564.99"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.14"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.51"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.378"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.84"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.36"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.16"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.22"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.05"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.77"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.83"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0531"
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.600"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0232"
$ws.Range("E51").Value = "  -3.50%  "
